$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Case's"
$ws.Range("B1").Value = "Operation"
$ws.Range("C1").Value = "Expexeped outcome"
$ws.Range("D1").Value = "True outcome"
$ws.Range("E1").Value = "Pass/Fail"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Clicking Start Button"
$ws.Range("C2").Value = "They game should move to the game screen "
$ws.Range("D2").Value = "Goes to game screen "
$ws.Range("E2").Value = "PASS"

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Clicking on rules button "
$ws.Range("C3").Value = "Should display a alert explaing the aleart"
$ws.Range("D3").Value = "Show user rules and how to play "
$ws.Range("E3").Value = "PASS"

# Column widths (character units, matches stored widths of
# 6.5 / 32.6640625 / 44.5 / 31.83203125 as closely as this engine allows)
$ws.Columns.Item(1).ColumnWidth = 5.666666666666667
$ws.Columns.Item(2).ColumnWidth = 31.830729166666668
$ws.Columns.Item(3).ColumnWidth = 43.666666666666664
$ws.Columns.Item(4).ColumnWidth = 30.998697916666668

# Selection matches the saved cursor position from the authored file
$ws.Range("E3").Select()
